# Update the cryptos list (price + 1h volume change) to the latest snapshot,
# and fix the swapped Polkadot / WrappedEther rows (13 <-> 14).
#
# Column D ("Price") holds plain text like "1.002" / "313.90", even though it
# looks numeric. If written with a plain .Value assignment while the cell is
# still General-formatted, Excel "smart" type-detection would silently turn
# it into a real number and drop significant trailing zeros (e.g. "313.90"
# -> 313.9, "1.002" -> 1.002 as a float, "0.07300" -> 0.073). To keep these
# as literal text (matching the source file's inlineStr cells), the whole
# Price column is temporarily switched to Text format ("@") before writing,
# then the formatting is cleared again afterwards so the cells end up with
# the same (default/no) style as before - only their text content changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$priceCol = $ws.Range("D2:D51")
$priceCol.NumberFormat = "@"

# --- Rows 13 & 14: swap Polkadot and WrappedEther, with their new values ---
$ws.Cells.Item(13, 2).Value = "WrappedEther"
$ws.Cells.Item(13, 3).Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Cells.Item(13, 4).Value = "1.860.11"
$ws.Cells.Item(13, 5).Value = "  +2.78%  "

$ws.Cells.Item(14, 2).Value = "Polkadot"
$ws.Cells.Item(14, 3).Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Cells.Item(14, 4).Value = "5.404"
$ws.Cells.Item(14, 5).Value = "  +1.27%  "

# --- Price (column D) and Volume(1h) (column E) refresh for every other row ---
$ws.Cells.Item(2, 4).Value = "27.249.98"
$ws.Cells.Item(2, 5).Value = "  +1.07%  "

$ws.Cells.Item(3, 4).Value = "1.856.72"
$ws.Cells.Item(3, 5).Value = "  +1.67%  "

$ws.Cells.Item(4, 4).Value = "1.002"
$ws.Cells.Item(4, 5).Value = "  -0.31%  "

$ws.Cells.Item(5, 4).Value = "313.90"
$ws.Cells.Item(5, 5).Value = "  +0.96%  "

$ws.Cells.Item(6, 4).Value = "1.001"
$ws.Cells.Item(6, 5).Value = "  -0.22%  "

$ws.Cells.Item(7, 4).Value = "0.4655"
$ws.Cells.Item(7, 5).Value = "  +0.41%  "

$ws.Cells.Item(8, 4).Value = "0.3710"
$ws.Cells.Item(8, 5).Value = "  +0.35%  "

$ws.Cells.Item(9, 4).Value = "0.07300"
$ws.Cells.Item(9, 5).Value = "  -0.61%  "

$ws.Cells.Item(10, 4).Value = "0.8920"
$ws.Cells.Item(10, 5).Value = "  +1.89%  "

$ws.Cells.Item(11, 4).Value = "20.08"
$ws.Cells.Item(11, 5).Value = "  +2.23%  "

$ws.Cells.Item(12, 4).Value = "0.07869"
$ws.Cells.Item(12, 5).Value = "  -0.03%  "

$ws.Cells.Item(15, 4).Value = "6.520"
$ws.Cells.Item(15, 5).Value = "  -0.39%  "

$ws.Cells.Item(16, 4).Value = "91.46"
$ws.Cells.Item(16, 5).Value = "  +0.06%  "

$ws.Cells.Item(17, 5).Value = "  -0.41%  "

$ws.Cells.Item(18, 4).Value = "0.000008929"
$ws.Cells.Item(18, 5).Value = "  +1.13%  "

$ws.Cells.Item(19, 5).Value = "  -0.20%  "

$ws.Cells.Item(20, 4).Value = "14.75"
$ws.Cells.Item(20, 5).Value = "  -0.09%  "

$ws.Cells.Item(21, 4).Value = "27.271.89"
$ws.Cells.Item(21, 5).Value = "  +1.06%  "

$ws.Cells.Item(22, 4).Value = "5.091"
$ws.Cells.Item(22, 5).Value = "  -0.18%  "

$ws.Cells.Item(23, 4).Value = "10.54"
$ws.Cells.Item(23, 5).Value = "  +0.16%  "

$ws.Cells.Item(24, 4).Value = "2.071.68"
$ws.Cells.Item(24, 5).Value = "  +2.34%  "

$ws.Cells.Item(25, 4).Value = "2.066"
$ws.Cells.Item(25, 5).Value = "  +11.42%  "

$ws.Cells.Item(26, 4).Value = "151.57"
$ws.Cells.Item(26, 5).Value = "  -0.31%  "

$ws.Cells.Item(27, 4).Value = "18.48"
$ws.Cells.Item(27, 5).Value = "  +0.27%  "

$ws.Cells.Item(28, 4).Value = "2.048"
$ws.Cells.Item(28, 5).Value = "  +0.43%  "

$ws.Cells.Item(29, 4).Value = "116.03"
$ws.Cells.Item(29, 5).Value = "  +0.37%  "

$ws.Cells.Item(30, 4).Value = "5.045"
$ws.Cells.Item(30, 5).Value = "  -1.18%  "

$ws.Cells.Item(31, 4).Value = "0.08819"
$ws.Cells.Item(31, 5).Value = "  -0.72%  "

$ws.Cells.Item(32, 4).Value = "3.143"
$ws.Cells.Item(32, 5).Value = "  +6.23%  "

$ws.Cells.Item(33, 4).Value = "0.7700"
$ws.Cells.Item(33, 5).Value = "  +5.33%  "

$ws.Cells.Item(34, 4).Value = "1.170"
$ws.Cells.Item(34, 5).Value = "  +3.45%  "

$ws.Cells.Item(35, 4).Value = "4.526"
$ws.Cells.Item(35, 5).Value = "  +1.97%  "

$ws.Cells.Item(36, 4).Value = "2.700"
$ws.Cells.Item(36, 5).Value = "  +9.36%  "

$ws.Cells.Item(37, 4).Value = "1.110"
$ws.Cells.Item(37, 5).Value = "  +3.33%  "

$ws.Cells.Item(38, 4).Value = "0.01944"
$ws.Cells.Item(38, 5).Value = "  +0.04%  "

$ws.Cells.Item(39, 4).Value = "0.05222"
$ws.Cells.Item(39, 5).Value = "  +0.12%  "

$ws.Cells.Item(40, 4).Value = "2.953"
$ws.Cells.Item(40, 5).Value = "  -0.19%  "

$ws.Cells.Item(41, 4).Value = "7.061"
$ws.Cells.Item(41, 5).Value = "  -0.48%  "

$ws.Cells.Item(42, 4).Value = "0.5114"
$ws.Cells.Item(42, 5).Value = "  -0.60%  "

$ws.Cells.Item(43, 4).Value = "0.1630"
$ws.Cells.Item(43, 5).Value = "  +0.17%  "

$ws.Cells.Item(44, 4).Value = "8.514"
$ws.Cells.Item(44, 5).Value = "  +4.42%  "

$ws.Cells.Item(45, 4).Value = "0.4803"
$ws.Cells.Item(45, 5).Value = "  -0.42%  "

$ws.Cells.Item(46, 5).Value = "  +1.61%  "

$ws.Cells.Item(47, 4).Value = "1.001"
$ws.Cells.Item(47, 5).Value = "  -0.27%  "

$ws.Cells.Item(48, 4).Value = "102.59"
$ws.Cells.Item(48, 5).Value = "  +0.81%  "

$ws.Cells.Item(49, 4).Value = "1.648"
$ws.Cells.Item(49, 5).Value = "  +1.54%  "

$ws.Cells.Item(50, 4).Value = "0.06198"
$ws.Cells.Item(50, 5).Value = "  -0.13%  "

$ws.Cells.Item(51, 5).Value = "  +1.27%  "

# Restore the column's original (default) styling now that every value is
# safely stored as text.
$priceCol.ClearFormats()
